# Update column G ("K") values on the active sheet to reflect the
# regenerated strikeout/mean/std stats (s_vals), per commit message:
# "regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = [ordered]@{
    2  = 8
    3  = 10
    4  = 8
    5  = 11
    6  = 8
    7  = 8
    8  = 6
    9  = 9
    10 = 5
    11 = 5
    12 = 3
    13 = 6
    14 = 7
    15 = 7
    16 = 8
    17 = 6
    18 = 6
    19 = 9
    20 = 6
    21 = 7
    22 = 6
    23 = 3
    24 = 6
    25 = 7
    26 = 4
    27 = 4
    28 = 1
    29 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
